$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gesamtinvestitionskosten")

# Update the base cost figures (column B, rows 2-10) that feed all the
# dependent formulas (C, D, E, F columns and the totals in rows 12/14,
# plus the cross-sheet reference on "Mittelverwendung - Mittelherkun").
$ws.Range("B2").Value = 110.0
$ws.Range("B3").Value = 1100.0
$ws.Range("B4").Value = 1100.0
$ws.Range("B5").Value = 110.0
$ws.Range("B6").Value = 110.0
$ws.Range("B7").Value = 110.0
$ws.Range("B8").Value = 110.0
$ws.Range("B9").Value = 110.0
$ws.Range("B10").Value = 110.0

# D10 is a hard-coded value (not a formula).
$ws.Range("D10").Value = 0.6

# Tax rate inputs used throughout the C/D column formulas.
$ws.Range("B20").Value = 0.5
$ws.Range("B21").Value = 0.33

$excel.CalculateFullRebuild()
